$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.769.71'
$ws.Range('E2').Value = '  +1.57%  '

$ws.Range('D3').Value = '1.720.31'
$ws.Range('E3').Value = '  +0.27%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9999'

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '239.99'
$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('E6').Value = '  +0.27%  '

$ws.Range('E7').Value = '  -2.25%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2555'
$ws.Range('E8').Value = '  -0.82%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06112'
$ws.Range('E9').Value = '  -0.90%  '

$ws.Range('D10').Value = '1.718.04'
$ws.Range('E10').Value = '  +0.25%  '

$ws.Range('E11').Value = '  +2.48%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06899'
$ws.Range('E12').Value = '  -0.66%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.5945'
$ws.Range('E13').Value = '  -0.37%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.402'
$ws.Range('E14').Value = '  -1.65%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '76.30'
$ws.Range('E15').Value = '  -0.14%  '

$ws.Range('D17').Value = '26.680.59'
$ws.Range('E17').Value = '  +1.29%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9998'
$ws.Range('E18').Value = '  +0.29%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007010'
$ws.Range('E19').Value = '  -1.08%  '

$ws.Range('E20').Value = '  -0.07%  '

$ws.Range('D21').Value = '1.938.99'
$ws.Range('E21').Value = '  +0.43%  '

$ws.Range('E22').Value = '  -1.05%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.305'
$ws.Range('E23').Value = '  -1.37%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.048'

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '140.83'
$ws.Range('E25').Value = '  +3.32%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.08'
$ws.Range('E26').Value = '  -0.47%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.777'
$ws.Range('E27').Value = '  +2.94%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '105.57'
$ws.Range('E28').Value = '  +0.07%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.368'
$ws.Range('E29').Value = '  -2.19%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.932'
$ws.Range('E30').Value = '  +1.65%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07870'
$ws.Range('E31').Value = '  -1.12%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.613'
$ws.Range('E32').Value = '  +0.07%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04545'
$ws.Range('E33').Value = '  +2.97%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.585'
$ws.Range('E34').Value = '  -0.51%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9888'
$ws.Range('E35').Value = '  -0.36%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6095'
$ws.Range('E36').Value = '  -1.12%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9236'
$ws.Range('E37').Value = '  -1.12%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.486'
$ws.Range('E38').Value = '  +4.79%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.958'
$ws.Range('E39').Value = '  -1.28%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9993'
$ws.Range('E40').Value = '  +0.27%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.711'
$ws.Range('E41').Value = '  +5.21%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '100.05'
$ws.Range('E43').Value = '  +0.47%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3779'
$ws.Range('E44').Value = '  -0.71%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.690'

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1140'
$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05349'
$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.720'
$ws.Range('E48').Value = '  +0.37%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.53'
$ws.Range('E49').Value = '  -2.94%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.226'
$ws.Range('E50').Value = '  +1.16%  '

$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.18%  '

